$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ligand/receptor expression statistics and derived edge weights
# with recomputed values based on updated TPM input data.
$ws.Cells.Item(2, 7).Value = 2.384498666666666
$ws.Cells.Item(2, 8).Value = 7.153495999999999
$ws.Cells.Item(2, 9).Value = 0.03997130847982828
$ws.Cells.Item(2, 10).Value = 0.03997130847982829
$ws.Cells.Item(2, 13).Value = 0.07919566666666666
$ws.Cells.Item(2, 14).Value = 0.237587
$ws.Cells.Item(2, 15).Value = 0.08232403487459106
$ws.Cells.Item(2, 16).Value = 0.08232403487459104
$ws.Cells.Item(2, 17).Value = 0.1888419615724444
$ws.Cells.Item(2, 18).Value = 1.699577654152
$ws.Cells.Item(2, 19).Value = 0.00329059939327642
$ws.Cells.Item(2, 20).Value = 0.00329059939327642
$ws.Cells.Item(3, 7).Value = 2.384498666666666
$ws.Cells.Item(3, 8).Value = 7.153495999999999
$ws.Cells.Item(3, 9).Value = 0.03997130847982828
$ws.Cells.Item(3, 10).Value = 0.03997130847982829
$ws.Cells.Item(3, 15).Value = 0.904479490283777
$ws.Cells.Item(3, 16).Value = 0.9044794902837769
$ws.Cells.Item(3, 17).Value = 2.074772955521777
$ws.Cells.Item(3, 18).Value = 18.67295659969599
$ws.Cells.Item(3, 19).Value = 0.0361532287198107
$ws.Cells.Item(3, 20).Value = 0.0361532287198107
$ws.Cells.Item(4, 7).Value = 2.384498666666666
$ws.Cells.Item(4, 8).Value = 7.153495999999999
$ws.Cells.Item(4, 9).Value = 0.03997130847982828
$ws.Cells.Item(4, 10).Value = 0.03997130847982829
$ws.Cells.Item(4, 13).Value = 0.012695
$ws.Cells.Item(4, 14).Value = 0.038085
$ws.Cells.Item(4, 15).Value = 0.01319647484163191
$ws.Cells.Item(4, 16).Value = 0.01319647484163191
$ws.Cells.Item(4, 17).Value = 0.03027121057333333
$ws.Cells.Item(4, 18).Value = 0.27244089516
$ws.Cells.Item(4, 19).Value = 0.0005274803667411621
$ws.Cells.Item(4, 20).Value = 0.0005274803667411622
$ws.Cells.Item(5, 9).Value = 0.9389652669334476
$ws.Cells.Item(5, 10).Value = 0.9389652669334477
$ws.Cells.Item(5, 13).Value = 0.07919566666666666
$ws.Cells.Item(5, 14).Value = 0.237587
$ws.Cells.Item(5, 15).Value = 0.08232403487459106
$ws.Cells.Item(5, 16).Value = 0.08232403487459104
$ws.Cells.Item(5, 17).Value = 4.436083020539334
$ws.Cells.Item(5, 18).Value = 39.924747184854
$ws.Cells.Item(5, 19).Value = 0.07729940938105884
$ws.Cells.Item(5, 20).Value = 0.07729940938105884
$ws.Cells.Item(6, 9).Value = 0.9389652669334476
$ws.Cells.Item(6, 10).Value = 0.9389652669334477
$ws.Cells.Item(6, 15).Value = 0.904479490283777
$ws.Cells.Item(6, 16).Value = 0.9044794902837769
$ws.Cells.Item(6, 19).Value = 0.8492748260301354
$ws.Cells.Item(6, 20).Value = 0.8492748260301354
$ws.Cells.Item(7, 9).Value = 0.9389652669334476
$ws.Cells.Item(7, 10).Value = 0.9389652669334477
$ws.Cells.Item(7, 13).Value = 0.012695
$ws.Cells.Item(7, 14).Value = 0.038085
$ws.Cells.Item(7, 15).Value = 0.01319647484163191
$ws.Cells.Item(7, 16).Value = 0.01319647484163191
$ws.Cells.Item(7, 17).Value = 0.7111004467300001
$ws.Cells.Item(7, 18).Value = 6.39990402057
$ws.Cells.Item(7, 19).Value = 0.01239103152225343
$ws.Cells.Item(7, 20).Value = 0.01239103152225343
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.036048
$ws.Cells.Item(8, 8).Value = 0.108144
$ws.Cells.Item(8, 9).Value = 0.0006042719789376481
$ws.Cells.Item(8, 10).Value = 0.0006042719789376483
$ws.Cells.Item(8, 13).Value = 0.07919566666666666
$ws.Cells.Item(8, 14).Value = 0.237587
$ws.Cells.Item(8, 15).Value = 0.08232403487459106
$ws.Cells.Item(8, 16).Value = 0.08232403487459104
$ws.Cells.Item(8, 17).Value = 0.002854845392
$ws.Cells.Item(8, 18).Value = 0.025693608528
$ws.Cells.Item(8, 19).Value = 0.0000497461074678011
$ws.Cells.Item(8, 20).Value = 0.0000497461074678011
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.036048
$ws.Cells.Item(9, 8).Value = 0.108144
$ws.Cells.Item(9, 9).Value = 0.0006042719789376481
$ws.Cells.Item(9, 10).Value = 0.0006042719789376483
$ws.Cells.Item(9, 15).Value = 0.904479490283777
$ws.Cells.Item(9, 16).Value = 0.9044794902837769
$ws.Cells.Item(9, 17).Value = 0.03136567721599999
$ws.Cells.Item(9, 18).Value = 0.282291094944
$ws.Cells.Item(9, 19).Value = 0.0005465516115022933
$ws.Cells.Item(9, 20).Value = 0.0005465516115022933
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.036048
$ws.Cells.Item(10, 8).Value = 0.108144
$ws.Cells.Item(10, 9).Value = 0.0006042719789376481
$ws.Cells.Item(10, 10).Value = 0.0006042719789376483
$ws.Cells.Item(10, 13).Value = 0.012695
$ws.Cells.Item(10, 14).Value = 0.038085
$ws.Cells.Item(10, 15).Value = 0.01319647484163191
$ws.Cells.Item(10, 16).Value = 0.01319647484163191
$ws.Cells.Item(10, 17).Value = 0.0004576293599999999
$ws.Cells.Item(10, 18).Value = 0.00411866424
$ws.Cells.Item(10, 19).Value = 0.0000079742599675538
$ws.Cells.Item(10, 20).Value = 0.000007974259967553801
$ws.Cells.Item(11, 7).Value = 1.188166
$ws.Cells.Item(11, 8).Value = 3.564498
$ws.Cells.Item(11, 9).Value = 0.01991720539631685
$ws.Cells.Item(11, 10).Value = 0.01991720539631685
$ws.Cells.Item(11, 13).Value = 0.07919566666666666
$ws.Cells.Item(11, 14).Value = 0.237587
$ws.Cells.Item(11, 15).Value = 0.08232403487459106
$ws.Cells.Item(11, 16).Value = 0.08232403487459104
$ws.Cells.Item(11, 17).Value = 0.09409759848066666
$ws.Cells.Item(11, 18).Value = 0.846878386326
$ws.Cells.Item(11, 19).Value = 0.001639664711650781
$ws.Cells.Item(11, 20).Value = 0.001639664711650781
$ws.Cells.Item(12, 7).Value = 1.188166
$ws.Cells.Item(12, 8).Value = 3.564498
$ws.Cells.Item(12, 9).Value = 0.01991720539631685
$ws.Cells.Item(12, 10).Value = 0.01991720539631685
$ws.Cells.Item(12, 15).Value = 0.904479490283777
$ws.Cells.Item(12, 16).Value = 0.9044794902837769
$ws.Cells.Item(12, 17).Value = 1.033833534038667
$ws.Cells.Item(12, 18).Value = 9.304501806347998
$ws.Cells.Item(12, 19).Value = 0.01801470378473796
$ws.Cells.Item(12, 20).Value = 0.01801470378473796
$ws.Cells.Item(13, 7).Value = 1.188166
$ws.Cells.Item(13, 8).Value = 3.564498
$ws.Cells.Item(13, 9).Value = 0.01991720539631685
$ws.Cells.Item(13, 10).Value = 0.01991720539631685
$ws.Cells.Item(13, 13).Value = 0.012695
$ws.Cells.Item(13, 14).Value = 0.038085
$ws.Cells.Item(13, 15).Value = 0.01319647484163191
$ws.Cells.Item(13, 16).Value = 0.01319647484163191
$ws.Cells.Item(13, 17).Value = 0.01508376737
$ws.Cells.Item(13, 18).Value = 0.13575390633
$ws.Cells.Item(13, 19).Value = 0.0002628368999281106
$ws.Cells.Item(13, 20).Value = 0.0002628368999281106
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.03233
$ws.Cells.Item(14, 8).Value = 0.09698999999999999
$ws.Cells.Item(14, 9).Value = 0.0005419472114695452
$ws.Cells.Item(14, 10).Value = 0.0005419472114695452
$ws.Cells.Item(14, 13).Value = 0.07919566666666666
$ws.Cells.Item(14, 14).Value = 0.237587
$ws.Cells.Item(14, 15).Value = 0.08232403487459106
$ws.Cells.Item(14, 16).Value = 0.08232403487459104
$ws.Cells.Item(14, 17).Value = 0.002560395903333333
$ws.Cells.Item(14, 18).Value = 0.02304356313
$ws.Cells.Item(14, 19).Value = 0.00004461528113720621
$ws.Cells.Item(14, 20).Value = 0.00004461528113720621
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.03233
$ws.Cells.Item(15, 8).Value = 0.09698999999999999
$ws.Cells.Item(15, 9).Value = 0.0005419472114695452
$ws.Cells.Item(15, 10).Value = 0.0005419472114695452
$ws.Cells.Item(15, 15).Value = 0.904479490283777
$ws.Cells.Item(15, 16).Value = 0.9044794902837769
$ws.Cells.Item(15, 17).Value = 0.02813061319333333
$ws.Cells.Item(15, 18).Value = 0.2531755187399999
$ws.Cells.Item(15, 19).Value = 0.0004901801375906886
$ws.Cells.Item(15, 20).Value = 0.0004901801375906885
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.03233
$ws.Cells.Item(16, 8).Value = 0.09698999999999999
$ws.Cells.Item(16, 9).Value = 0.0005419472114695452
$ws.Cells.Item(16, 10).Value = 0.0005419472114695452
$ws.Cells.Item(16, 13).Value = 0.012695
$ws.Cells.Item(16, 14).Value = 0.038085
$ws.Cells.Item(16, 15).Value = 0.01319647484163191
$ws.Cells.Item(16, 16).Value = 0.01319647484163191
$ws.Cells.Item(16, 17).Value = 0.00041042935
$ws.Cells.Item(16, 18).Value = 0.00369386415
$ws.Cells.Item(16, 19).Value = 0.00000715179274165042
$ws.Cells.Item(16, 20).Value = 0.00000715179274165042
